# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.124.28"
$ws.Range("E2").Value = "  +2.19%  "

# Row 3
$ws.Range("D3").Value = "3.464.43"
$ws.Range("E3").Value = "  +1.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7
$ws.Range("D7").Value = "3.465.10"
$ws.Range("E7").Value = "  +1.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.85%  "

# Row 13
$ws.Range("D13").Value = "4.057.17"
$ws.Range("E13").Value = "  +1.81%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.74%  "

# Row 15
$ws.Range("E15").Value = "  +2.63%  "

# Row 16
$ws.Range("D16").Value = "3.459.62"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17
$ws.Range("E17").Value = "  +0.85%  "

# Row 18
$ws.Range("D18").Value = "63.084.73"
$ws.Range("E18").Value = "  +2.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("E23").Value = "  +1.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "

# Row 25
$ws.Range("E25").Value = "  -0.19%  "

# Row 26
$ws.Range("D26").Value = "3.607.52"
$ws.Range("E26").Value = "  +1.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "

# Row 28
$ws.Range("E28").Value = "  -2.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.13%  "

# Row 32
$ws.Range("E32").Value = "  -1.26%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.24%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.29%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.60%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "170.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "

# Row 40
$ws.Range("E40").Value = "  +5.34%  "

# Row 41
$ws.Range("D41").Value = "3.501.38"
$ws.Range("E41").Value = "  +1.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0761"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.797"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
$ws.Range("E45").Value = "  +4.37%  "

# Row 46
$ws.Range("E46").Value = "  +2.35%  "

# Row 47
$ws.Range("E47").Value = "  -0.67%  "

# Row 48
$ws.Range("D48").Value = "2.617.99"
$ws.Range("E48").Value = "  +4.95%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "

